$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A5").Value = "AP1880.90E9.BE80"
$ws.Range("C5").Value = "172.19.4.145"
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 30
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 2
$ws.Range("A6").Select()
$ws.Range("A6").Value = ""
